$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.322.13"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.926.87"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.32"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.52"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.500"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.90"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.42"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.411.17"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.299.65"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.932.00"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "431.42"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.51"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  -1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.84"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  -2.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.72"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -4.98%  "
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("E30").Value = "  -2.86%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.62"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.109"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0879"
$ws.Range("E34").Value = "  +2.93%  "
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.62"
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("E37").Value = "  -2.53%  "
$ws.Range("E38").Value = "  -0.21%  "
$ws.Range("E39").Value = "  -2.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.55"
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("E41").Value = "  +4.44%  "
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.699.64"
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "365.77"
$ws.Range("E45").Value = "  -2.80%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.60"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("E48").Value = "  -2.35%  "
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.00"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("E51").Value = "  -1.36%  "
